$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 21, shifting existing rows 21-38 down to 22-39.
$ws.Rows.Item(21).Insert()

# Populate the newly inserted row 21 with the new weekly price record.
$ws.Cells.Item(21, 1).Value = 10
$ws.Cells.Item(21, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(21, 3).Value = "La Araucanía"
$ws.Cells.Item(21, 4).Value = 44452
$ws.Cells.Item(21, 5).Value = 9
$ws.Cells.Item(21, 6).Value = "Fruta"
$ws.Cells.Item(21, 7).Value = 100108
$ws.Cells.Item(21, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(21, 9).Value = 100108007
$ws.Cells.Item(21, 10).Value = "Coco"
$ws.Cells.Item(21, 11).Value = "Sin especificar"
$ws.Cells.Item(21, 12).Value = "Primera"
$ws.Cells.Item(21, 13).Value = 25
$ws.Cells.Item(21, 14).Value = 25000
$ws.Cells.Item(21, 15).Value = 25000
$ws.Cells.Item(21, 16).Value = 25000
$ws.Cells.Item(21, 17).Value = "$/malla 20 unidades"
$ws.Cells.Item(21, 18).Value = "Perú"
$ws.Cells.Item(21, 19).Value = 1250
$ws.Cells.Item(21, 20).Value = 20
